$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the stop name "Overtoom" to "Overtoom/ 1e C. Huygensstraat" (row 17, column A)
$ws.Range("A17").Value = "Overtoom/ 1e C. Huygensstraat"

# Reflect the resulting UI selection state from the edit
$ws.Range("F13").Select()
